# "adresse RZO.xlsx" - update subnetting example from 172.16.0.0/26 (mis-labelled /58)
# to 172.16.1.0/26, and rework the "ASYMETRIQUE" VLSM table (rows 11-14) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Header CIDR note (F1) ---
$ws.Range("F1").Value = "CIDR = 32 - 6 =26"

# --- Symmetric table (rows 3-6) and the matching cells in the asymmetric
#     table (rows 11-14) that keep the same addresses ---
$ws.Range("B3").Value = "172.16.1.0/26"
$ws.Range("B11").Value = "172.16.1.0/26"

$ws.Range("B4").Value = "172.16.1.64/26"
$ws.Range("B5").Value = "172.16.1.128/26"
$ws.Range("B6").Value = "172.16.1.192/26"

$ws.Range("C3").Value = "172.16.1.63"
$ws.Range("C11").Value = "172.16.1.63"

$ws.Range("C4").Value = "172.16.1.127"
$ws.Range("C13").Value = "172.16.1.127"

$ws.Range("C5").Value = "172.16.1.191"
$ws.Range("C6").Value = "172.16.1.255"

$ws.Range("D3").Value = "172.16.1.1"
$ws.Range("D11").Value = "172.16.1.1"

$ws.Range("D4").Value = "172.16.1.65"
$ws.Range("D12").Value = "172.16.1.65"

$ws.Range("D5").Value = "172.16.1.129"
$ws.Range("D14").Value = "172.16.1.129"

$ws.Range("D6").Value = "172.16.1.193"

$ws.Range("E3").Value = "172.16.1.62"
$ws.Range("E11").Value = "172.16.1.62"

$ws.Range("E4").Value = "172.16.1.126"
$ws.Range("E13").Value = "172.16.1.126"

$ws.Range("E5").Value = "172.16.1.190"
$ws.Range("E6").Value = "172.16.1.254"

# --- Asymmetric / VLSM table new addresses (rows 11-14) ---
$ws.Range("B12").Value = "172.16.1.64/27"
$ws.Range("B13").Value = "172.16.1.96/27"
$ws.Range("B14").Value = "172.16.1.128/28"

$ws.Range("C12").Value = "172.16.1.95"
$ws.Range("C14").Value = "172.16.1.143"

$ws.Range("D13").Value = "172.16.1.97"

$ws.Range("E12").Value = "172.16.1.94"
$ws.Range("E14").Value = "172.16.1.142"

# --- Update the selected range shown when the sheet is next opened ---
[void]$ws.Range("A9:E14").Select()
